{"js": "// The document has a bulleted \"Consolidated:\" list near the end of the\n// document that repeats (and starts with) \"Search patron records\" before\n// \"Check out books\". The edit removes that leading \"Search patron\n// records\" bullet item (its paragraph is deleted), which shifts every\n// following bullet's text up by one position and shortens the list by a\n// single item. Word also relocates its automatic \"_GoBack\" bookmark (the\n// marker for \"last edit location\") from the end of the list to the start\n// of the new first bullet (\"Check out books\").\n\nconst body = context.document.body;\n\n// Find every paragraph whose text is exactly \"Search patron records\".\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// There are a few \"Search patron records\" bullets earlier in the document\n// (one per user-type section); the one we need to remove is the last\n// occurrence, which lives in the final \"Consolidated:\" list.\nlet targetIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"Search patron records\") {\n    targetIndex = i;\n  }\n}\n\nif (targetIndex === -1) {\n  throw new Error('Could not find a paragraph with text \"Search patron records\"');\n}\n\n// Remove that paragraph entirely (this is what merges/shifts the rest of\n// the list up by one entry).\nparagraphs.items[targetIndex].delete();\nawait context.sync();\n\n// Move the \"_GoBack\" bookmark: drop it from wherever it currently sits and\n// re-insert it at the very start of the paragraph that is now first in\n// line (this is the paragraph that used to be second, holding \"Check out\n// books\", and it now occupies the slot the deleted paragraph vacated).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst refreshedParagraphs = body.paragraphs;\nrefreshedParagraphs.load(\"items/text\");\nawait context.sync();\n\nconst newFirstListItem = refreshedParagraphs.items[targetIndex];\nconst startOfItem = newFirstListItem.getRange(\"Start\");\nstartOfItem.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# The document has a bulleted \"Consolidated:\" list near the end of the\n# document that repeats (and starts with) \"Search patron records\" before\n# \"Check out books\". The edit removes that leading \"Search patron\n# records\" bullet item (its paragraph is deleted), which shifts every\n# following bullet's text up by one position and shortens the list by a\n# single item. Word also relocates its automatic \"_GoBack\" bookmark (the\n# marker for \"last edit location\") from the end of the list to the start\n# of the new first bullet (\"Check out books\").\n\n$d = $word.ActiveDocument\n\n# Find the LAST paragraph whose text equals \"Search patron records\" -- the\n# document repeats this bullet once per user-type section; the occurrence\n# we must remove is the final one, inside the \"Consolidated:\" list.\n$targetIndex = -1\n$i = 0\nforeach ($p in $d.Paragraphs) {\n    $i = $i + 1\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq \"Search patron records\") {\n        $targetIndex = $i\n    }\n}\n\nif ($targetIndex -eq -1) {\n    throw \"Could not find a paragraph with text 'Search patron records'\"\n}\n\n# Remove that whole bullet paragraph; this merges it away and shifts every\n# following bullet's text up by one position.\n$target = $d.Paragraphs.Item($targetIndex)\n$target.Range.Delete()\n\n# Word's automatic \"_GoBack\" bookmark (the \"last edit location\" marker)\n# moves to the edit point: drop the old one...\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# ...and re-create it, collapsed, at the very start of the paragraph that\n# is now first in line (previously second, holding \"Check out books\").\n$newFirst = $d.Paragraphs.Item($targetIndex)\n$startRange = $d.Range($newFirst.Range.Start, $newFirst.Range.Start)\n$d.Bookmarks.Add(\"_GoBack\", $startRange)\n"}
